# Scheduled runner update: refresh computed market-price / profit figures
# on several leve rows across the ALC, ARM, BSM, CRP, GSM, LTW and WVR
# sheets of the Pandaemonium_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 258.75
$ws.Range("I33").Value = 248.6842
$ws.Range("K33").Value = 248.6842
$ws.Range("M33").Value = -19.6842
$ws.Range("H41").Value = 100
$ws.Range("I41").Value = 100
$ws.Range("K41").Value = 100
$ws.Range("M41").Value = 340
$ws.Range("H70").Value = 1383.381
$ws.Range("I70").Value = 1555.6
$ws.Range("J70").Value = 1329.5625
$ws.Range("K70").Value = 4666.799999999999
$ws.Range("L70").Value = 3988.6875
$ws.Range("M70").Value = -4396.799999999999
$ws.Range("N70").Value = -4528.6875
$ws.Range("H73").Value = 1383.381
$ws.Range("I73").Value = 1555.6
$ws.Range("J73").Value = 1329.5625
$ws.Range("K73").Value = 4666.799999999999
$ws.Range("L73").Value = 3988.6875
$ws.Range("M73").Value = -3730.799999999999
$ws.Range("N73").Value = -5860.6875
$ws.Range("H137").Value = 1923.371
$ws.Range("I137").Value = 1482.5581
$ws.Range("J137").Value = 2921
$ws.Range("K137").Value = 4447.6743
$ws.Range("L137").Value = 8763
$ws.Range("M137").Value = -1897.6743
$ws.Range("N137").Value = -13863
$ws.Range("H138").Value = 4118297
$ws.Range("I138").Value = 1681.4482
$ws.Range("J138").Value = 6414102
$ws.Range("K138").Value = 5044.3446
$ws.Range("L138").Value = 19242306
$ws.Range("M138").Value = 95.65539999999964
$ws.Range("N138").Value = -19252586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17558.37
$ws.Range("I32").Value = 18868.908
$ws.Range("K32").Value = 18868.908
$ws.Range("M32").Value = -18581.908
$ws.Range("H61").Value = 11157.645
$ws.Range("I61").Value = 8516.429
$ws.Range("J61").Value = 20401.9
$ws.Range("K61").Value = 8516.429
$ws.Range("L61").Value = 20401.9
$ws.Range("M61").Value = -8304.429
$ws.Range("N61").Value = -20825.9
$ws.Range("H63").Value = 3815.8333
$ws.Range("I63").Value = 3581
$ws.Range("J63").Value = 4990
$ws.Range("K63").Value = 3581
$ws.Range("L63").Value = 4990
$ws.Range("M63").Value = -2895
$ws.Range("N63").Value = -6362
$ws.Range("H66").Value = 3815.8333
$ws.Range("I66").Value = 3581
$ws.Range("J66").Value = 4990
$ws.Range("K66").Value = 17905
$ws.Range("L66").Value = 24950
$ws.Range("M66").Value = -14473
$ws.Range("N66").Value = -31814
$ws.Range("H105").Value = 45000
$ws.Range("J105").Value = 45000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -51988
$ws.Range("H136").Value = 11157.645
$ws.Range("I136").Value = 8516.429
$ws.Range("J136").Value = 20401.9
$ws.Range("K136").Value = 25549.287
$ws.Range("L136").Value = 61205.7
$ws.Range("M136").Value = -22999.287
$ws.Range("N136").Value = -66305.70000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1702.1666
$ws.Range("I99").Value = 1610.8182
$ws.Range("J99").Value = 1845.7142
$ws.Range("K99").Value = 1610.8182
$ws.Range("L99").Value = 1845.7142
$ws.Range("M99").Value = -112.8181999999999
$ws.Range("N99").Value = -4841.7142
$ws.Range("H106").Value = 20985.572
$ws.Range("J106").Value = 20985.572
$ws.Range("L106").Value = 20985.572
$ws.Range("N106").Value = -23509.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6091.1934
$ws.Range("I31").Value = 8325.1875
$ws.Range("J31").Value = 3708.2666
$ws.Range("K31").Value = 8325.1875
$ws.Range("L31").Value = 3708.2666
$ws.Range("M31").Value = -8030.1875
$ws.Range("N31").Value = -4298.2666
$ws.Range("H34").Value = 6091.1934
$ws.Range("I34").Value = 8325.1875
$ws.Range("J34").Value = 3708.2666
$ws.Range("K34").Value = 8325.1875
$ws.Range("L34").Value = 3708.2666
$ws.Range("M34").Value = -8123.1875
$ws.Range("N34").Value = -4112.2666
$ws.Range("H48").Value = 14012.25
$ws.Range("J48").Value = 14012.25
$ws.Range("L48").Value = 14012.25
$ws.Range("N48").Value = -14964.25
$ws.Range("H58").Value = 1717557.8
$ws.Range("I58").Value = 2526798.2
$ws.Range("J58").Value = 3872.1177
$ws.Range("K58").Value = 2526798.2
$ws.Range("L58").Value = 3872.1177
$ws.Range("M58").Value = -2526595.2
$ws.Range("N58").Value = -4278.1177
$ws.Range("H105").Value = 829.9
$ws.Range("I105").Value = 590.2308
$ws.Range("J105").Value = 2387.75
$ws.Range("K105").Value = 590.2308
$ws.Range("L105").Value = 2387.75
$ws.Range("M105").Value = 1156.7692
$ws.Range("N105").Value = -5881.75
$ws.Range("H132").Value = 4643.512
$ws.Range("I132").Value = 5013.7856
$ws.Range("J132").Value = 3846
$ws.Range("K132").Value = 15041.3568
$ws.Range("L132").Value = 11538
$ws.Range("M132").Value = -12511.3568
$ws.Range("N132").Value = -16598
$ws.Range("H136").Value = 1717557.8
$ws.Range("I136").Value = 2526798.2
$ws.Range("J136").Value = 3872.1177
$ws.Range("K136").Value = 7580394.600000001
$ws.Range("L136").Value = 11616.3531
$ws.Range("M136").Value = -7577844.600000001
$ws.Range("N136").Value = -16716.3531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1247
$ws.Range("I43").Value = 1247
$ws.Range("K43").Value = 1247
$ws.Range("M43").Value = -1096
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6798.6665
$ws.Range("I122").Value = 6430.091
$ws.Range("K122").Value = 19290.273
$ws.Range("M122").Value = -16840.273
$ws.Range("H132").Value = 4789.7104
$ws.Range("I132").Value = 4626.161
$ws.Range("K132").Value = 13878.483
$ws.Range("M132").Value = -11348.483

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 665
$ws.Range("I100").Value = 629
$ws.Range("J100").Value = 701
$ws.Range("K100").Value = 1258
$ws.Range("L100").Value = 1402
$ws.Range("M100").Value = -717
$ws.Range("N100").Value = -2484
$ws.Range("H132").Value = 1372.1608
$ws.Range("I132").Value = 548.4
$ws.Range("J132").Value = 2745.0952
$ws.Range("K132").Value = 1645.2
$ws.Range("L132").Value = 8235.285600000001
$ws.Range("M132").Value = 884.8000000000002
$ws.Range("N132").Value = -13295.2856
